$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "value" figures (province fund data) that must be rescaled
# (multiplied by 10000) to push/pull the data onto the correct magnitude.
# Row 29's D cell is intentionally blank and is left untouched.
# Target values are assigned explicitly (rather than via in-script
# multiplication) so the stored double bit pattern exactly matches the
# scaled figures.
$targets = @{
    2  = 1904669.259845
    3  = 34538854.574683
    4  = 521448.463218
    5  = 1296957.950658
    6  = 6625774.632704
    7  = 73131.686734
    8  = 259626.007438
    9  = 21860.72564
    10 = 395437.700184
    11 = 450442.369135
    12 = 1272637.271388
    13 = 7893406.600703
    14 = 1050577.584345
    15 = 264697.101112
    16 = 5729178.227088
    17 = 648483.799995
    18 = 2674926.938553
    19 = 153871.5165
    20 = 35201.810472
    21 = 468848.702909
    22 = 98357896.370719
    23 = 3075545.091631
    24 = 512007.135779
    25 = 7218211.048453
    26 = 6637837.9737
    27 = 2823201.624189
    28 = 5034368.645675
    30 = 30321.574091
    31 = 518364.337086
    32 = 980965.540184
    33 = 572366.30633
}

foreach ($row in $targets.Keys) {
    $ws.Cells.Item($row, 4).Value2 = $targets[$row]
}
